# The workbook's data dictionary (row 4 of Sheet1) documents the schema's
# scalar type names. This rename swaps the old "number"/"string" type
# spellings for "float"/"str" everywhere they appear in the header labels
# (and the nested container/tuple descriptions that mention them).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "float id"
$ws.Range("C4").Value = "float n"
$ws.Range("D4").Value = "str s"
$ws.Range("F4").Value = "[float] n_list"
$ws.Range("G4").Value = "[str] s_list"
$ws.Range("I4").Value = "{float} n_dict"
$ws.Range("J4").Value = "{str} s_dict"
$ws.Range("K4").Value = "<bool b, float n, strs> t"
$ws.Range("L4").Value = "[<bool b, float n, str s>] t_list"
$ws.Range("M4").Value = "[{float}] d_list"
$ws.Range("N4").Value = "{<bool b, float n, str s>} t_dict"
$ws.Range("O4").Value = "<<bool b, float n, str s> t> t_type"

# Move the active selection to where the editor left off.
$ws.Activate()
$ws.Range("M6").Select()
